# Chart示例.xlsx — "A new version, which looks pretty fine as a demo..."
#
# Walks through the same sequence of sheet visits / edits the author made:
#   - on "taps", tweak a stray point value
#   - on "slides", drop a duplicated sample row and rescroll/reselect
#   - on "stars", rescale a batch of tangent offsets (F/H columns) and
#     finish there (so "stars" ends up the active tab)
#   - "planes" loses its old selection/active-tab state along the way
#
# Final active sheet/tab == "stars" (workbook bookViews.activeTab == 5).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "taps" sheet (sheet2.xml): single value tweak + navigate/select
# ---------------------------------------------------------------------
$taps = $wb.Worksheets.Item("taps")
$taps.Activate()
$taps.Range("B616").Value = -1.2
$taps.Range("F620").Select()
$excel.ActiveWindow.ScrollRow = 140
$excel.ActiveWindow.ScrollColumn = 1

# ---------------------------------------------------------------------
# 2) "slides" sheet (sheet4.xml): delete the duplicated row 192
#    (116.336 / 0 / 1.2 / 1), which shifts every following row up by
#    one and shrinks the used range from A1:D317 to A1:D316.
# ---------------------------------------------------------------------
$slides = $wb.Worksheets.Item("slides")
$slides.Activate()
$slides.Rows.Item(192).Delete()
$slides.Range("J197").Select()
$excel.ActiveWindow.ScrollRow = 182
$excel.ActiveWindow.ScrollColumn = 1

# ---------------------------------------------------------------------
# 3) "stars" sheet (sheet6.xml): rescale a batch of F/H offsets
# ---------------------------------------------------------------------
$stars = $wb.Worksheets.Item("stars")
$stars.Activate()

$stars.Range("F14").Value = -1
$stars.Range("H15").Value = -1
$stars.Range("F16").Value = -1
$stars.Range("F17").Value = 1
$stars.Range("H18").Value = 1

$stars.Range("H39").Value = -1
$stars.Range("F40").Value = -1
$stars.Range("H41").Value = 1
$stars.Range("F42").Value = 1
$stars.Range("H43").Value = -1
$stars.Range("F44").Value = -1

$stars.Range("H46").Value = -1
$stars.Range("F47").Value = -1
$stars.Range("H48").Value = 1
$stars.Range("F49").Value = 1

$stars.Range("H56").Value = 1
$stars.Range("F57").Value = 1
$stars.Range("H58").Value = -1
$stars.Range("F59").Value = -1

$stars.Range("F75").Value = -1.2
$stars.Range("H76").Value = 0.6
$stars.Range("F77").Value = 0.6
$stars.Range("H78").Value = 1.2
$stars.Range("F79").Value = 1.2
$stars.Range("H80").Value = -0.6
$stars.Range("F81").Value = -0.6
$stars.Range("H82").Value = -1.2

$stars.Range("H23").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1

# ---------------------------------------------------------------------
# 4) "planes" sheet (sheet1.xml): move selection, leave it not-active
#    (tabSelected only remains on whichever sheet is active at save
#    time — "stars" — so visiting planes last-but-one clears its flag)
# ---------------------------------------------------------------------
$planes = $wb.Worksheets.Item("planes")
$planes.Activate()
$planes.Range("C31").Select()

# ---------------------------------------------------------------------
# End on "stars" so it is the active tab / sheet on save.
# ---------------------------------------------------------------------
$stars.Activate()
$stars.Range("H23").Select()
